# Auto-generated script applying the 2024-12-04 crime data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 7408
$ws.Range('K3').Value = 7675
$ws.Range('C4').Value = 1854
$ws.Range('J4').Value = 1846
$ws.Range('K4').Value = 1614
$ws.Range('K5').Value = 545
$ws.Range('K6').Value = 8504
$ws.Range('C7').Value = 28399
$ws.Range('J7').Value = 29315
$ws.Range('K7').Value = 25746

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K2').Value = 84
$ws.Range('K7').Value = 316

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 465
$ws.Range('K6').Value = 559
$ws.Range('K7').Value = 1672

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K2').Value = 271
$ws.Range('K3').Value = 387
$ws.Range('K6').Value = 346
$ws.Range('K7').Value = 1087

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 245
$ws.Range('K3').Value = 283
$ws.Range('K6').Value = 262
$ws.Range('K7').Value = 861

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K6').Value = 226
$ws.Range('K7').Value = 608

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K2').Value = 115
$ws.Range('K3').Value = 184
$ws.Range('K7').Value = 436

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 223
$ws.Range('K7').Value = 778
$ws.Range('K8').Value = 1672
$ws.Range('K9').Value = 122
$ws.Range('K15').Value = 263
$ws.Range('K18').Value = 172
$ws.Range('K19').Value = 744
$ws.Range('K20').Value = 634
$ws.Range('K23').Value = 259
$ws.Range('K25').Value = 119
$ws.Range('K29').Value = 1417
$ws.Range('K32').Value = 28
$ws.Range('K33').Value = 1087
$ws.Range('K34').Value = 148
$ws.Range('K37').Value = 861
$ws.Range('K42').Value = 946
$ws.Range('K47').Value = 176
$ws.Range('K49').Value = 144
$ws.Range('K50').Value = 118
$ws.Range('K51').Value = 324
$ws.Range('K53').Value = 316
$ws.Range('K54').Value = 506
$ws.Range('C63').Value = 281
$ws.Range('J63').Value = 126
$ws.Range('K63').Value = 74
$ws.Range('K65').Value = 608
$ws.Range('K67').Value = 1005
$ws.Range('K75').Value = 81
$ws.Range('K78').Value = 317
$ws.Range('K79').Value = 631
$ws.Range('K84').Value = 207
$ws.Range('K85').Value = 1180
$ws.Range('K99').Value = 436
$ws.Range('C101').Value = 28399
$ws.Range('J101').Value = 29315
$ws.Range('K101').Value = 25746

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 279
$ws.Range('K3').Value = 363
$ws.Range('K7').Value = 1005

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K2').Value = 67
$ws.Range('K6').Value = 40
$ws.Range('K7').Value = 207

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K6').Value = 71
$ws.Range('K7').Value = 144

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K3').Value = 115
$ws.Range('K6').Value = 279
$ws.Range('K7').Value = 506

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 402
$ws.Range('K5').Value = 35
$ws.Range('K7').Value = 1417

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 218
$ws.Range('K4').Value = 34
$ws.Range('K6').Value = 247
$ws.Range('K7').Value = 744

$ws = $wb.Worksheets.Item('River North')
$ws.Range('K3').Value = 69
$ws.Range('K6').Value = 177

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 258
$ws.Range('K3').Value = 280
$ws.Range('K7').Value = 946

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K3').Value = 82
$ws.Range('K6').Value = 106
$ws.Range('K7').Value = 317

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('K2').Value = 71
$ws.Range('K7').Value = 259

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K3').Value = 201
$ws.Range('K6').Value = 163
$ws.Range('K7').Value = 631

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 216
$ws.Range('K6').Value = 181
$ws.Range('K7').Value = 634

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K3').Value = 60
$ws.Range('K7').Value = 172

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K4').Value = 30
$ws.Range('K6').Value = 216
$ws.Range('K7').Value = 778

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K6').Value = 42
$ws.Range('K7').Value = 148

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K2').Value = 47
$ws.Range('K7').Value = 119

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K2').Value = 50
$ws.Range('K7').Value = 176

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K6').Value = 77
$ws.Range('K7').Value = 263

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K4').Value = 14
$ws.Range('K7').Value = 118

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('K3').Value = 45
$ws.Range('K7').Value = 122

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K2').Value = 67
$ws.Range('K6').Value = 69
$ws.Range('K7').Value = 223

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range('K4').Value = 3
$ws.Range('K7').Value = 28

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('K2').Value = 29
$ws.Range('K7').Value = 81

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K6').Value = 104
$ws.Range('K7').Value = 324

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 388
$ws.Range('K3').Value = 406
$ws.Range('K7').Value = 1180
